$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) edits to remain plain text, matching the source
# workbook where these cells are stored as inline/shared strings rather than
# numbers (values use "." as a thousands separator, e.g. "25.864.81").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.864.81"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.619.99"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.72"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.46"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.38"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.56"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.874.83"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.49"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.50"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.89"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.65"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.19"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126.11"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.839"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.511"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.30"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.755.51"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0113"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.01"

# Restore the default (unstyled) cell style now that the text is safely stored,
# so no extra formatting is introduced versus the original workbook.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"

# Column E (Volume(1h)) values are already safe as text (leading/trailing
# spaces + "%" sign), so a direct assignment is sufficient.

$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -1.32%  "
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("E18").Value = "  -3.62%  "
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("E22").Value = "  -3.22%  "
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -7.37%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  -4.54%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -6.80%  "
$ws.Range("E44").Value = "  -5.91%  "
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  -4.67%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("E51").Value = "  -1.04%  "
